# "molar vol included in gsa now"
#
# Adds two new result columns -- "Molar Vol (L/mol)" and
# "M. Vol. err (L/mol)" -- to the per-step results table on the
# "Gas Sorption Input" sheet, between the existing "Mass frac./Mass
# frac. err" columns (M) and the "Dual Mode Pred" columns (old N:Q,
# now shifted to P:S).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gas Sorption Input")

# Insert two blank columns at N:O, pushing the old N:Q (Mass frac.,
# Mass frac. err, Dual Mode Pred (CC/CC), Dual Mode Pred Err (CC/CC))
# out to P:S. This also fixes up the sheet dimension / row spans
# automatically.
$ws.Columns("N:O").Insert()

# New header row (row 11) labels for the inserted columns.
$ws.Range("N11").Value = "Molar Vol (L/mol)"
$ws.Range("O11").Value = "M. Vol. err (L/mol)"

# New per-step data values (rows 12-18).
$ws.Range("N12").Value = 11.437863057602693
$ws.Range("O12").Value = 0.005789915463944882

$ws.Range("N13").Value = 3.430803399487487
$ws.Range("O13").Value = 0.0017860378557054417

$ws.Range("N14").Value = 2.0046748445944713
$ws.Range("O14").Value = 0.0010746059098845988

$ws.Range("N15").Value = 1.3310716164493712
$ws.Range("O15").Value = 0.0007402721728057062

$ws.Range("N16").Value = 0.9847644037805116
$ws.Range("O16").Value = 0.0005700445884418672

$ws.Range("N17").Value = 0.7710205285078638
$ws.Range("O17").Value = 0.0004666441566139297

$ws.Range("N18").Value = 0.6458420268502967
$ws.Range("O18").Value = 0.00040748136774608885
